$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Row 38 already had its date (C38 = 10/14/2014) and the "started at
# midnight" note (D38) logged. Record what that session's work actually
# was and how many hours it took.
$ws.Range("A38").Value = "Writing Simple Ai Script"
$ws.Range("B38").Value = 0.5

# A second entry was logged the same day (row 39): recording the script
# that had just been written. Copy row 38's date cell (so the new date
# cell keeps the same date number format/style) then set its value.
$ws.Range("C38").Copy() | Out-Null
$ws.Range("C39").PasteSpecial(-4122) | Out-Null

$ws.Range("A39").Value = "Recording Simple Ai Script"
$ws.Range("B39").Value = 1.5
$ws.Range("C39").Value = 41926

# Leave the sheet scrolled/selected the same way the author left it.
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H28").Select()
